$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..7) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
